$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 93: the timestamp in A93 was wrong (it had a time-of-day
# component instead of lining up with the 07:00 snapshot the rest of the
# series uses). Correct it in place; everything else in the row is right.
$ws.Cells.Item(93, 1).Value = 45476.2916666667

# --- Row 94: newly pulled bar appended by the refreshed R script.
$ws.Cells.Item(94, 1).Value = 45477.6495949074
$ws.Cells.Item(93, 1).Copy()
$ws.Cells.Item(94, 1).PasteSpecial(-4122)   # xlPasteFormats - match the date style used by the rest of column A

$ws.Cells.Item(94, 2).Value = 7500
$ws.Cells.Item(94, 3).Value = 3.38000011444092
$ws.Cells.Item(94, 4).Value = 3.29999995231628
$ws.Cells.Item(94, 5).Value = 3.29999995231628
$ws.Cells.Item(94, 6).Value = 3.36999988555908

# adj_close is stored as text in this sheet (mirrors the 'close' value).
# A leading apostrophe forces text entry instead of Excel auto-converting
# it back to a number, then we drop the resulting quote-prefix formatting
# so the cell carries no explicit style, matching the rest of the column.
$ws.Cells.Item(94, 7).Value = "'3.36999988555908"
$ws.Cells.Item(94, 7).Style = "Normal"

$ws.Cells.Item(94, 8).Value = "ESPE.MI"

$excel.Application.CutCopyMode = $false
